$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell text updates from the diff.
# D-column price cells get NumberFormat "@" forced before the write so
# numeric-looking text (e.g. "2.41") is not silently coerced to a number,
# then the style is reset to "Normal" afterwards so no new number format
# sticks to the cell (matches original un-styled price/volume cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.341.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.799.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("E6").Value = "  +3.83%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.293"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0688"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0963"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.059.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.824.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.637"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.335.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0783"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0518"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("E34").Value = "  -3.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.373.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.654"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.21%  "
$ws.Range("E37").Value = "  -1.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.56%  "
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.944"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.960.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.53%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("E51").Value = "  +1.37%  "
